{"js": "// Bold the six \"section heading\" lines that introduce each SQL query block.\n// (The rest of the source diff is Word's own proofing-tool markup\n// (w:proofErr / run-splitting) added automatically when the file was\n// re-opened and spell/grammar-checked in the desktop app; it carries no\n// textual or formatting change and is not something this script needs to\n// reproduce \u2014 the only author-visible edit is the new bold emphasis.)\n\nconst headings = [\n  \"SUM (TOTAL AMOUNT PAID TO EMPLOYEE AS SALARY):\",\n  \"COUNT+ GROUP_BY  (NUMBER OF EMPLOYEES IN EACH DEPARTMENT):\",\n  \"AVG+ GROUP_BY (AVERAGE OF SALARY IN EACH DEPARTMENT):\",\n  \"MIN, MAX (MINIMUM AND MAXIMUM SALARY IN EACH DEPARTMENT):\",\n  \"SUM (TOTAL EXPERIENCE OF EMPLOYEE IN EACH DEPARTMENT):\",\n  \"GROUP_BY+ HAVING (EMPLOYEE MORE THAN ONE YEAR OF EXPERIENCE):\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (headings.indexOf(para.text) !== -1) {\n    para.font.bold = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Bold the six \"section heading\" lines that introduce each SQL query block.\n# (The rest of the source diff is Word's own proofing-tool markup\n# (w:proofErr / run-splitting) added automatically when the file was\n# re-opened and spell/grammar-checked in the desktop app; it carries no\n# textual or formatting change and is not something this script needs to\n# reproduce -- the only author-visible edit is the new bold emphasis.)\n\n$d = $word.ActiveDocument\n\n$headings = @(\n    \"SUM (TOTAL AMOUNT PAID TO EMPLOYEE AS SALARY):\",\n    \"COUNT+ GROUP_BY  (NUMBER OF EMPLOYEES IN EACH DEPARTMENT):\",\n    \"AVG+ GROUP_BY (AVERAGE OF SALARY IN EACH DEPARTMENT):\",\n    \"MIN, MAX (MINIMUM AND MAXIMUM SALARY IN EACH DEPARTMENT):\",\n    \"SUM (TOTAL EXPERIENCE OF EMPLOYEE IN EACH DEPARTMENT):\",\n    \"GROUP_BY+ HAVING (EMPLOYEE MORE THAN ONE YEAR OF EXPERIENCE):\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`n\", \"`a\", [char]7)\n    if ($headings -contains $text) {\n        $p.Range.Font.Bold = 1\n    }\n}\n"}
